$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the now-unused "issue_code" column (column B) entirely.
$ws.Columns.Item(2).Delete()

# Swap the order of the two "Golden Parachute" rows (old rows 4 and 5).
$ws.Range("A4").Value = "Approval of Golden Parachute"
$ws.Range("A5").Value = "Approval of Golden Parachute Payments"

# Append two new proposal rows at the bottom of the list.
$ws.Range("A28").Value = "Appointment of Samwise as auditor for fiscal year 2024"
$ws.Range("A29").Value = "Director election of silly sally"

# Match the resulting view/selection state.
$ws.Range("A30").Select()
